$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Enter when 7 day cases >"
$ws.Range("A13").Value = "Leave when 7 day cases <"

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 5600
$ws.Range("D12").Value = 11200
$ws.Range("E12").Value = 22400
$ws.Range("F12").Value = "For moderate supress and moderate_suppress_no_4, the only policy in the model"

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 3733
$ws.Range("D13").Value = 7466
$ws.Range("E13").Value = 14933
$ws.Range("F13").Value = "For moderate supress and moderate_suppress_no_4, the only policy in the model"

$ws.Range("F16").Select()
